$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("873:873").Insert()

$ws.Range("A873").Value = 3
$ws.Range("B873").Value = "Femacal de La Calera"
$ws.Range("C873").Value = "Coquimbo"
$ws.Range("D873").Value = 45132
$ws.Range("E873").Value = 5
$ws.Range("F873").Value = 100112045
$ws.Range("G873").Value = "Zapallo"
$ws.Range("H873").Value = "Camote"
$ws.Range("I873").Value = "1a (guarda)"
$ws.Range("J873").Value = 210
$ws.Range("K873").Value = 480
$ws.Range("L873").Value = 500
$ws.Range("M873").Value = 490
$ws.Range("N873").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O873").Value = "Provincia de Talca"
$ws.Range("P873").Value = 490
$ws.Range("Q873").Value = 1
$ws.Range("R873").Value = "Hortaliza"
